$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-10 05:18:15"
$ws.Range("E3").Value = "2026-02-10 05:18:18"
$ws.Range("I3").Value = "6.6 mm"
$ws.Range("O3").Value = "-0.7 °C"
$ws.Range("E4").Value = "2026-02-10 05:18:20"
$ws.Range("J4").Value = "1004.9 hPa"
$ws.Range("K4").Value = "-0.1 MJ/m2"
$ws.Range("N4").Value = "7.2 °C 4:52 TU"
$ws.Range("O4").Value = "8.9 °C"
$ws.Range("E5").Value = "2026-02-10 05:18:22"
$ws.Range("G5").Value = "130 cm"
$ws.Range("I5").Value = "10.7 mm"
$ws.Range("E6").Value = "2026-02-10 05:18:25"
$ws.Range("J6").Value = "1005.0 hPa"
$ws.Range("N6").Value = "6.6 °C 4:48 TU"
$ws.Range("O6").Value = "7.4 °C"
$ws.Range("E7").Value = "2026-02-10 05:18:27"
$ws.Range("E8").Value = "2026-02-10 05:18:29"
$ws.Range("L8").Value = "42.8 km/h - 311º 4:42 TU"
$ws.Range("E9").Value = "2026-02-10 05:18:32"
$ws.Range("N9").Value = "5.1 °C 4:45 TU"
$ws.Range("O9").Value = "6.7 °C"
$ws.Range("E10").Value = "2026-02-10 05:18:34"
$ws.Range("N10").Value = "6.3 °C 4:59 TU"
$ws.Range("O10").Value = "7.3 °C"
$ws.Range("E11").Value = "2026-02-10 05:18:37"
$ws.Range("E12").Value = "2026-02-10 05:18:39"
$ws.Range("N12").Value = "5.5 °C 4:50 TU"
$ws.Range("O12").Value = "7.0 °C"
$ws.Range("E13").Value = "2026-02-10 05:18:41"
$ws.Range("H13").Value = "'96%"
$ws.Range("I13").Value = "1.8 mm"
$ws.Range("J13").Value = "1008.1 hPa"
$ws.Range("N13").Value = "2.4 °C 4:45 TU"
$ws.Range("E14").Value = "2026-02-10 05:18:44"
$ws.Range("N14").Value = "7.7 °C 4:59 TU"
$ws.Range("O14").Value = "9.4 °C"
$ws.Range("E15").Value = "2026-02-10 05:18:46"
$ws.Range("N15").Value = "3.9 °C 4:53 TU"
$ws.Range("O15").Value = "6.6 °C"
$ws.Range("E16").Value = "2026-02-10 05:18:48"
$ws.Range("G16").Value = "77 cm"
$ws.Range("I16").Value = "10.1 mm"
$ws.Range("O16").Value = "-0.4 °C"
$ws.Range("E17").Value = "2026-02-10 05:18:51"
$ws.Range("H17").Value = "'94%"
$ws.Range("M17").Value = "5.4 °C 4:57 TU"
$ws.Range("O17").Value = "2.3 °C"
$ws.Range("E18").Value = "2026-02-10 05:18:53"
$ws.Range("H18").Value = "'98%"
$ws.Range("N18").Value = "6.1 °C 4:59 TU"
$ws.Range("O18").Value = "7.7 °C"
$ws.Range("E19").Value = "2026-02-10 05:18:56"
$ws.Range("E20").Value = "2026-02-10 05:18:58"
$ws.Range("O20").Value = "-1.2 °C"
$ws.Range("E21").Value = "2026-02-10 05:19:00"
$ws.Range("I21").Value = "2.7 mm"
$ws.Range("L21").Value = "8.6 km/h - 30º 4:35 TU"
$ws.Range("O21").Value = "4.3 °C"
$ws.Range("E22").Value = "2026-02-10 05:19:03"
$ws.Range("H22").Value = "'99%"
$ws.Range("I22").Value = "0.1 mm"
$ws.Range("E23").Value = "2026-02-10 05:19:05"
$ws.Range("I23").Value = "8.2 mm"
$ws.Range("M23").Value = "0.4 °C 4:39 TU"
$ws.Range("O23").Value = "-0.4 °C"
$ws.Range("E24").Value = "2026-02-10 05:19:07"
$ws.Range("E25").Value = "2026-02-10 05:19:10"
$ws.Range("H25").Value = "'94%"
$ws.Range("I25").Value = "5.4 mm"
$ws.Range("E26").Value = "2026-02-10 05:19:12"
$ws.Range("M26").Value = "3.5 °C 4:59 TU"
$ws.Range("O26").Value = "2.9 °C"
$ws.Range("E27").Value = "2026-02-10 05:19:15"
$ws.Range("I27").Value = "1.5 mm"
$ws.Range("E28").Value = "2026-02-10 05:19:17"
$ws.Range("I28").Value = "0.1 mm"
$ws.Range("J28").Value = "1005.5 hPa"
$ws.Range("N28").Value = "3.9 °C 4:40 TU"
$ws.Range("O28").Value = "5.2 °C"
$ws.Range("E29").Value = "2026-02-10 05:19:19"
$ws.Range("H29").Value = "'92%"
$ws.Range("N29").Value = "7.3 °C 4:39 TU"
$ws.Range("O29").Value = "9.1 °C"
$ws.Range("E30").Value = "2026-02-10 05:19:22"
$ws.Range("L30").Value = "16.6 km/h - 327º 4:59 TU"
$ws.Range("N30").Value = "6.9 °C 4:39 TU"
$ws.Range("E31").Value = "2026-02-10 05:19:24"
$ws.Range("H31").Value = "'86%"
$ws.Range("O31").Value = "8.9 °C"
$ws.Range("E32").Value = "2026-02-10 05:19:27"
$ws.Range("E33").Value = "2026-02-10 05:19:29"
$ws.Range("I33").Value = "2.8 mm"
$ws.Range("E34").Value = "2026-02-10 05:19:32"
$ws.Range("I34").Value = "1.7 mm"
$ws.Range("M34").Value = "3.7 °C 4:38 TU"
$ws.Range("O34").Value = "2.6 °C"
$ws.Range("E35").Value = "2026-02-10 05:19:34"
$ws.Range("I35").Value = "0.1 mm"
$ws.Range("J35").Value = "1005.5 hPa"
$ws.Range("N35").Value = "9.7 °C 4:56 TU"
$ws.Range("E36").Value = "2026-02-10 05:19:37"
$ws.Range("N36").Value = "6.9 °C 4:50 TU"
$ws.Range("O36").Value = "9.1 °C"
$ws.Range("E37").Value = "2026-02-10 05:19:39"
$ws.Range("J37").Value = "1007.0 hPa"
$ws.Range("O37").Value = "3.9 °C"
$ws.Range("E38").Value = "2026-02-10 05:19:42"
$ws.Range("N38").Value = "6.9 °C 4:56 TU"
$ws.Range("O38").Value = "7.8 °C"
$ws.Range("E39").Value = "2026-02-10 05:19:44"
$ws.Range("I39").Value = "1.7 mm"
$ws.Range("L39").Value = "48.6 km/h - 346º 4:42 TU"
$ws.Range("E40").Value = "2026-02-10 05:19:46"
$ws.Range("I40").Value = "3.2 mm"
$ws.Range("J40").Value = "1008.2 hPa"
$ws.Range("N40").Value = "4.2 °C 4:53 TU"
$ws.Range("E41").Value = "2026-02-10 05:19:49"
$ws.Range("N41").Value = "8.2 °C 4:34 TU"
$ws.Range("O41").Value = "9.9 °C"
$ws.Range("E42").Value = "2026-02-10 05:19:51"
$ws.Range("N42").Value = "7.1 °C 4:35 TU"
$ws.Range("O42").Value = "8.2 °C"
$ws.Range("E43").Value = "2026-02-10 05:19:53"
$ws.Range("N43").Value = "5.4 °C 4:57 TU"
$ws.Range("E44").Value = "2026-02-10 05:19:56"
$ws.Range("I44").Value = "6.8 mm"
$ws.Range("E45").Value = "2026-02-10 05:19:58"
$ws.Range("I45").Value = "13.3 mm"
$ws.Range("E46").Value = "2026-02-10 05:20:00"
$ws.Range("H46").Value = "'100%"
$ws.Range("J46").Value = "1006.7 hPa"
